$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 776.12
$ws.Range("I33").Value = 598.7368
$ws.Range("J33").Value = 1337.8334
$ws.Range("K33").Value = 598.7368
$ws.Range("L33").Value = 1337.8334
$ws.Range("M33").Value = -369.7368
$ws.Range("N33").Value = -1795.8334
# Row 40
$ws.Range("H40").Value = 4532.5713
$ws.Range("J40").Value = 4796.615
$ws.Range("L40").Value = 4796.615
$ws.Range("N40").Value = -5146.615
# Row 106
$ws.Range("H106").Value = 4544.2666
$ws.Range("I106").Value = 3826.75
$ws.Range("J106").Value = 5364.2856
$ws.Range("K106").Value = 3826.75
$ws.Range("L106").Value = 5364.2856
$ws.Range("M106").Value = -3195.75
$ws.Range("N106").Value = -6626.2856
# Row 113
$ws.Range("H113").Value = 3964.5417
$ws.Range("I113").Value = 3314.9
$ws.Range("J113").Value = 4428.5713
$ws.Range("K113").Value = 3314.9
$ws.Range("L113").Value = 4428.5713
$ws.Range("M113").Value = -60.90000000000009
$ws.Range("N113").Value = -10936.5713
# Row 116
$ws.Range("H116").Value = 56207.023
$ws.Range("I116").Value = 87268.60000000001
$ws.Range("J116").Value = 4437.7334
$ws.Range("K116").Value = 87268.60000000001
$ws.Range("L116").Value = 4437.7334
$ws.Range("M116").Value = -83826.60000000001
$ws.Range("N116").Value = -11321.7334
# Row 132
$ws.Range("H132").Value = 4068.8235
$ws.Range("I132").Value = 2415.487
$ws.Range("J132").Value = 9442.166999999999
$ws.Range("K132").Value = 7246.461
$ws.Range("L132").Value = 28326.501
$ws.Range("M132").Value = -4716.461
$ws.Range("N132").Value = -33386.501
# Row 137
$ws.Range("H137").Value = 3345.9744
$ws.Range("I137").Value = 3627.5715
$ws.Range("K137").Value = 10882.7145
$ws.Range("M137").Value = -8332.7145
# Row 138
$ws.Range("H138").Value = 2791.2534
$ws.Range("J138").Value = 3068.9106
$ws.Range("L138").Value = 9206.731800000001
$ws.Range("N138").Value = -19486.7318

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1661.55
$ws.Range("I2").Value = 2021.5
$ws.Range("J2").Value = 821.6667
$ws.Range("K2").Value = 2021.5
$ws.Range("L2").Value = 821.6667
$ws.Range("M2").Value = -1908.5
$ws.Range("N2").Value = -1047.6667
# Row 32
$ws.Range("H32").Value = 23080.482
$ws.Range("I32").Value = 8991.73
$ws.Range("J32").Value = 44012.344
$ws.Range("K32").Value = 8991.73
$ws.Range("L32").Value = 44012.344
$ws.Range("M32").Value = -8704.73
$ws.Range("N32").Value = -44586.344
# Row 116
$ws.Range("H116").Value = 1661.55
$ws.Range("I116").Value = 2021.5
$ws.Range("J116").Value = 821.6667
$ws.Range("K116").Value = 2021.5
$ws.Range("L116").Value = 821.6667
$ws.Range("M116").Value = 272.5
$ws.Range("N116").Value = -5409.6667

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1661.55
$ws.Range("I3").Value = 2021.5
$ws.Range("J3").Value = 821.6667
$ws.Range("K3").Value = 2021.5
$ws.Range("L3").Value = 821.6667
$ws.Range("M3").Value = -1907.5
$ws.Range("N3").Value = -1049.6667
# Row 80
$ws.Range("H80").Value = 460.72223
$ws.Range("I80").Value = 219.75
$ws.Range("J80").Value = 529.5714
$ws.Range("K80").Value = 219.75
$ws.Range("L80").Value = 529.5714
$ws.Range("M80").Value = 778.25
$ws.Range("N80").Value = -2525.5714
# Row 83
$ws.Range("H83").Value = 460.72223
$ws.Range("I83").Value = 219.75
$ws.Range("J83").Value = 529.5714
$ws.Range("K83").Value = 1098.75
$ws.Range("L83").Value = 2647.857
$ws.Range("M83").Value = 3893.25
$ws.Range("N83").Value = -12631.857
# Row 86
$ws.Range("H86").Value = 3840.739
$ws.Range("I86").Value = 3112.5454
$ws.Range("J86").Value = 5689.231
$ws.Range("K86").Value = 3112.5454
$ws.Range("L86").Value = 5689.231
$ws.Range("M86").Value = -1989.5454
$ws.Range("N86").Value = -7935.231
# Row 89
$ws.Range("H89").Value = 3840.739
$ws.Range("I89").Value = 3112.5454
$ws.Range("J89").Value = 5689.231
$ws.Range("K89").Value = 15562.727
$ws.Range("L89").Value = 28446.155
$ws.Range("M89").Value = -9946.726999999999
$ws.Range("N89").Value = -39678.155
# Row 94
$ws.Range("H94").Value = 1584.0476
$ws.Range("I94").Value = 904.7857
$ws.Range("J94").Value = 2942.5715
$ws.Range("K94").Value = 904.7857
$ws.Range("L94").Value = 2942.5715
$ws.Range("M94").Value = -453.7857
$ws.Range("N94").Value = -3844.5715
# Row 99
$ws.Range("H99").Value = 3656.0557
$ws.Range("I99").Value = 4712.154
$ws.Range("J99").Value = 910.2
$ws.Range("K99").Value = 4712.154
$ws.Range("L99").Value = 910.2
$ws.Range("M99").Value = -3214.154
$ws.Range("N99").Value = -3906.2
# Row 134
$ws.Range("H134").Value = 30416.432
$ws.Range("I134").Value = 44478.61
$ws.Range("J134").Value = 7314.2856
$ws.Range("K134").Value = 133435.83
$ws.Range("L134").Value = 21942.8568
$ws.Range("M134").Value = -130900.83
$ws.Range("N134").Value = -27012.8568

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2385.747
$ws.Range("I31").Value = 1711.2656
$ws.Range("J31").Value = 4262.5654
$ws.Range("K31").Value = 1711.2656
$ws.Range("L31").Value = 4262.5654
$ws.Range("M31").Value = -1416.2656
$ws.Range("N31").Value = -4852.5654
# Row 34
$ws.Range("H34").Value = 2385.747
$ws.Range("I34").Value = 1711.2656
$ws.Range("J34").Value = 4262.5654
$ws.Range("K34").Value = 1711.2656
$ws.Range("L34").Value = 4262.5654
$ws.Range("M34").Value = -1509.2656
$ws.Range("N34").Value = -4666.5654
# Row 105
$ws.Range("H105").Value = 596.52856
$ws.Range("I105").Value = 592.3333
$ws.Range("J105").Value = 634.2857
$ws.Range("K105").Value = 592.3333
$ws.Range("L105").Value = 634.2857
$ws.Range("M105").Value = 1154.6667
$ws.Range("N105").Value = -4128.2857
# Row 107
$ws.Range("H107").Value = 428.33334
$ws.Range("I107").Value = 655
$ws.Range("J107").Value = 383
$ws.Range("K107").Value = 655
$ws.Range("L107").Value = 383
$ws.Range("M107").Value = 1265
$ws.Range("N107").Value = -4223
# Row 122
$ws.Range("H122").Value = 2466
$ws.Range("I122").Value = 2819.875
$ws.Range("J122").Value = 1657.1428
$ws.Range("K122").Value = 8459.625
$ws.Range("L122").Value = 4971.428400000001
$ws.Range("M122").Value = -6009.625
$ws.Range("N122").Value = -9871.428400000001
# Row 132
$ws.Range("H132").Value = 1925.2115
$ws.Range("I132").Value = 866.9429
$ws.Range("J132").Value = 4104
$ws.Range("K132").Value = 2600.8287
$ws.Range("L132").Value = 12312
$ws.Range("M132").Value = -70.82870000000003
$ws.Range("N132").Value = -17372
# Row 134
$ws.Range("H134").Value = 1581.7
$ws.Range("I134").Value = 1004.4483
$ws.Range("J134").Value = 2378.8572
$ws.Range("K134").Value = 3013.3449
$ws.Range("L134").Value = 7136.571599999999
$ws.Range("M134").Value = -478.3449000000001
$ws.Range("N134").Value = -12206.5716

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 98
$ws.Range("H98").Value = 9549
$ws.Range("J98").Value = 14173.5
$ws.Range("L98").Value = 42520.5
$ws.Range("N98").Value = -45516.5
# Row 134
$ws.Range("H134").Value = 4649.5
$ws.Range("I134").Value = 4597.1113
$ws.Range("J134").Value = 4806.6665
$ws.Range("K134").Value = 13791.3339
$ws.Range("L134").Value = 14419.9995
$ws.Range("M134").Value = -8721.333899999998
$ws.Range("N134").Value = -24559.9995

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 6892.647
$ws.Range("I113").Value = 9067.5
$ws.Range("K113").Value = 9067.5
$ws.Range("M113").Value = -6897.5
# Row 126
$ws.Range("H126").Value = 3468.182
$ws.Range("I126").Value = 4699.75
$ws.Range("J126").Value = 2764.4285
$ws.Range("K126").Value = 14099.25
$ws.Range("L126").Value = 8293.2855
$ws.Range("M126").Value = -11629.25
$ws.Range("N126").Value = -13233.2855
# Row 132
$ws.Range("H132").Value = 3161.1853
$ws.Range("I132").Value = 3038.8125
$ws.Range("J132").Value = 3339.182
$ws.Range("K132").Value = 9116.4375
$ws.Range("L132").Value = 10017.546
$ws.Range("M132").Value = -6586.4375
$ws.Range("N132").Value = -15077.546

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 477.72726
$ws.Range("I22").Value = 437.83334
$ws.Range("J22").Value = 525.6
$ws.Range("K22").Value = 437.83334
$ws.Range("L22").Value = 525.6
$ws.Range("M22").Value = -142.83334
$ws.Range("N22").Value = -1115.6
# Row 27
$ws.Range("H27").Value = 477.72726
$ws.Range("I27").Value = 437.83334
$ws.Range("J27").Value = 525.6
$ws.Range("K27").Value = 437.83334
$ws.Range("L27").Value = 525.6
$ws.Range("M27").Value = -330.83334
$ws.Range("N27").Value = -739.6
# Row 122
$ws.Range("H122").Value = 3736.2727
$ws.Range("I122").Value = 4340
$ws.Range("J122").Value = 3233.1667
$ws.Range("K122").Value = 13020
$ws.Range("L122").Value = 9699.500100000001
$ws.Range("M122").Value = -10570
$ws.Range("N122").Value = -14599.5001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 51982.125
$ws.Range("J46").Value = 51982.125
$ws.Range("L46").Value = 51982.125
$ws.Range("N46").Value = -52444.125
# Row 113
$ws.Range("H113").Value = 333.42307
$ws.Range("I113").Value = 261.26666
$ws.Range("J113").Value = 431.81818
$ws.Range("K113").Value = 783.79998
$ws.Range("L113").Value = 1295.45454
$ws.Range("M113").Value = 1386.20002
$ws.Range("N113").Value = -5635.45454
# Row 134
$ws.Range("H134").Value = 51982.125
$ws.Range("J134").Value = 51982.125
$ws.Range("L134").Value = 155946.375
$ws.Range("N134").Value = -161016.375
# Row 136
$ws.Range("H136").Value = 381877.44
$ws.Range("I136").Value = 73245.42999999999
$ws.Range("J136").Value = 669934
$ws.Range("K136").Value = 219736.29
$ws.Range("L136").Value = 2009802
$ws.Range("M136").Value = -217186.29
$ws.Range("N136").Value = -2014902
